# Auto-generated from diff: updates cached market/profit values across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1335.1
$ws.Range("I28").Value = 1335.1
$ws.Range("K28").Value = 1335.1
$ws.Range("M28").Value = -850.0999999999999

$ws.Range("H53").Value = 1631.8889
$ws.Range("I53").Value = 537.5
$ws.Range("K53").Value = 537.5
$ws.Range("M53").Value = 99.5

$ws.Range("H62").Value = 9590.409
$ws.Range("I62").Value = 9020.0625
$ws.Range("K62").Value = 9020.0625
$ws.Range("M62").Value = -8396.0625

$ws.Range("H65").Value = 9590.409
$ws.Range("I65").Value = 9020.0625
$ws.Range("K65").Value = 45100.3125
$ws.Range("M65").Value = -41980.3125

$ws.Range("H80").Value = 296611.2
$ws.Range("I80").Value = 12840
$ws.Range("J80").Value = 1053334.4
$ws.Range("K80").Value = 38520
$ws.Range("L80").Value = 3160003.2
$ws.Range("M80").Value = -37522
$ws.Range("N80").Value = -3161999.2

$ws.Range("H83").Value = 296611.2
$ws.Range("I83").Value = 12840
$ws.Range("J83").Value = 1053334.4
$ws.Range("K83").Value = 115560
$ws.Range("L83").Value = 9480009.6
$ws.Range("M83").Value = -110568
$ws.Range("N83").Value = -9489993.6

$ws.Range("H92").Value = 3856
$ws.Range("I92").Value = 759.4286
$ws.Range("K92").Value = 759.4286
$ws.Range("M92").Value = 488.5714

$ws.Range("H106").Value = 900
$ws.Range("I106").Value = 900
$ws.Range("K106").Value = 900
$ws.Range("M106").Value = -269

$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1300
$ws.Range("K113").Value = 1300
$ws.Range("M113").Value = 1954

$ws.Range("H132").Value = 11868.4
$ws.Range("I132").Value = 12158.762
$ws.Range("K132").Value = 36476.286
$ws.Range("M132").Value = -33946.286

$ws.Range("H137").Value = 8105.132
$ws.Range("I137").Value = 13779.186
$ws.Range("J137").Value = 2212.8462
$ws.Range("K137").Value = 41337.558
$ws.Range("L137").Value = 6638.5386
$ws.Range("M137").Value = -38787.558
$ws.Range("N137").Value = -11738.5386

$ws.Range("H138").Value = 3923.1016
$ws.Range("J138").Value = 4294.2095
$ws.Range("L138").Value = 12882.6285
$ws.Range("N138").Value = -23162.6285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6640.747
$ws.Range("I32").Value = 6297.026
$ws.Range("J32").Value = 12002.8
$ws.Range("K32").Value = 6297.026
$ws.Range("L32").Value = 12002.8
$ws.Range("M32").Value = -6010.026
$ws.Range("N32").Value = -12576.8

$ws.Range("H45").Value = 4400.528
$ws.Range("I45").Value = 3816.6924
$ws.Range("K45").Value = 3816.6924
$ws.Range("M45").Value = -3439.6924

$ws.Range("H61").Value = 4117.6387
$ws.Range("I61").Value = 3341.4
$ws.Range("K61").Value = 3341.4
$ws.Range("M61").Value = -3129.4

$ws.Range("H74").Value = 548562.5600000001
$ws.Range("I74").Value = 1001115.5
$ws.Range("K74").Value = 1001115.5
$ws.Range("M74").Value = -1000241.5

$ws.Range("H77").Value = 548562.5600000001
$ws.Range("I77").Value = 1001115.5
$ws.Range("K77").Value = 5005577.5
$ws.Range("M77").Value = -5001209.5

$ws.Range("H97").Value = 1326.5834
$ws.Range("I97").Value = 975.45
$ws.Range("K97").Value = 975.45
$ws.Range("M97").Value = -479.45

$ws.Range("H102").Value = 1513.1842
$ws.Range("I102").Value = 1369
$ws.Range("J102").Value = 3195.3333
$ws.Range("K102").Value = 1369
$ws.Range("L102").Value = 3195.3333
$ws.Range("M102").Value = 253
$ws.Range("N102").Value = -6439.3333

$ws.Range("H132").Value = 2501.7073
$ws.Range("I132").Value = 2139.2334
$ws.Range("J132").Value = 3490.2727
$ws.Range("K132").Value = 6417.7002
$ws.Range("L132").Value = 10470.8181
$ws.Range("M132").Value = -3887.7002
$ws.Range("N132").Value = -15530.8181

$ws.Range("H136").Value = 4117.6387
$ws.Range("I136").Value = 3341.4
$ws.Range("K136").Value = 10024.2
$ws.Range("M136").Value = -7474.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2749.6667
$ws.Range("I99").Value = 2530.75
$ws.Range("J99").Value = 3187.5
$ws.Range("K99").Value = 2530.75
$ws.Range("L99").Value = 3187.5
$ws.Range("M99").Value = -1032.75
$ws.Range("N99").Value = -6183.5

$ws.Range("H107").Value = 89751.336
$ws.Range("J107").Value = 9998.5
$ws.Range("L107").Value = 9998.5
$ws.Range("N107").Value = -13838.5

$ws.Range("H134").Value = 1672.3572
$ws.Range("I134").Value = 1414.3195
$ws.Range("K134").Value = 4242.958500000001
$ws.Range("M134").Value = -1707.958500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 95
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 95
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -321

$ws.Range("H31").Value = 5722.393
$ws.Range("I31").Value = 4420.25
$ws.Range("J31").Value = 7458.5835
$ws.Range("K31").Value = 4420.25
$ws.Range("L31").Value = 7458.5835
$ws.Range("M31").Value = -4125.25
$ws.Range("N31").Value = -8048.5835

$ws.Range("H34").Value = 5722.393
$ws.Range("I34").Value = 4420.25
$ws.Range("J34").Value = 7458.5835
$ws.Range("K34").Value = 4420.25
$ws.Range("L34").Value = 7458.5835
$ws.Range("M34").Value = -4218.25
$ws.Range("N34").Value = -7862.5835

$ws.Range("H94").Value = 3428.6667
$ws.Range("J94").Value = 3996
$ws.Range("L94").Value = 3996
$ws.Range("N94").Value = -4898

$ws.Range("H120").Value = 20719.7
$ws.Range("J120").Value = 21149.625
$ws.Range("L120").Value = 21149.625
$ws.Range("N120").Value = -28407.625

$ws.Range("H122").Value = 7275.222
$ws.Range("I122").Value = 6415.3335
$ws.Range("K122").Value = 19246.0005
$ws.Range("M122").Value = -16796.0005

$ws.Range("H124").Value = 25326
$ws.Range("J124").Value = 25326
$ws.Range("L124").Value = 25326
$ws.Range("N124").Value = -30236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2568.1765
$ws.Range("J113").Value = 2281.9333
$ws.Range("L113").Value = 6845.7999
$ws.Range("N113").Value = -11185.7999

$ws.Range("H131").Value = 144023.47
$ws.Range("I131").Value = 850953
$ws.Range("J131").Value = 2637.56
$ws.Range("K131").Value = 2552859
$ws.Range("L131").Value = 7912.68
$ws.Range("M131").Value = -2547819
$ws.Range("N131").Value = -17992.68

$ws.Range("H139").Value = 3925.7856
$ws.Range("I139").Value = 3175.3635
$ws.Range("K139").Value = 9526.0905
$ws.Range("M139").Value = -4386.0905

$ws.Range("H140").Value = 2539.5667
$ws.Range("I140").Value = 2539.5667
$ws.Range("K140").Value = 7618.7001
$ws.Range("M140").Value = -2438.7001

$ws.Range("H141").Value = 5179.231
$ws.Range("I141").Value = 4883.8696
$ws.Range("K141").Value = 14651.6088
$ws.Range("M141").Value = -9471.6088

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 52510.25
$ws.Range("J46").Value = 46666.668
$ws.Range("L46").Value = 46666.668
$ws.Range("N46").Value = -46978.668

$ws.Range("H101").Value = 21200
$ws.Range("J101").Value = 21200
$ws.Range("L101").Value = 21200
$ws.Range("N101").Value = -27690

$ws.Range("H107").Value = 387.25
$ws.Range("I107").Value = 387.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 387.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1532.75
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 3545.6667
$ws.Range("I113").Value = 2939.125
$ws.Range("K113").Value = 2939.125
$ws.Range("M113").Value = -769.125

$ws.Range("H132").Value = 4188.478
$ws.Range("J132").Value = 5039.5557
$ws.Range("L132").Value = 15118.6671
$ws.Range("N132").Value = -20178.6671

$ws.Range("H141").Value = 107777.2
$ws.Range("J141").Value = 107777.2
$ws.Range("L141").Value = 107777.2
$ws.Range("N141").Value = -118137.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 41016.176
$ws.Range("I25").Value = 23030.555
$ws.Range("K25").Value = 23030.555
$ws.Range("M25").Value = -22800.555

$ws.Range("H40").Value = 2784.8572
$ws.Range("I40").Value = 2582.3333
$ws.Range("K40").Value = 2582.3333
$ws.Range("M40").Value = -2446.3333

$ws.Range("H93").Value = 3187.4707
$ws.Range("I93").Value = 3091.2
$ws.Range("J93").Value = 3325
$ws.Range("K93").Value = 3091.2
$ws.Range("L93").Value = 3325
$ws.Range("M93").Value = -1843.2
$ws.Range("N93").Value = -5821

$ws.Range("H100").Value = 2941.6155
$ws.Range("I100").Value = 2924.3
$ws.Range("J100").Value = 2999.3333
$ws.Range("K100").Value = 2924.3
$ws.Range("L100").Value = 2999.3333
$ws.Range("M100").Value = -2383.3
$ws.Range("N100").Value = -4081.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 666.625
$ws.Range("I100").Value = 654.7143
$ws.Range("J100").Value = 750
$ws.Range("K100").Value = 1309.4286
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -768.4286
$ws.Range("N100").Value = -2582

$ws.Range("H132").Value = 6253221.5
$ws.Range("I132").Value = 8931916
$ws.Range("K132").Value = 26795748
$ws.Range("M132").Value = -26793218

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
